$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 169
$ws.Range("F5").Value = 1340
$ws.Range("F6").Value = 240
$ws.Range("F7").Value = 2555
$ws.Range("F8").Value = 942
$ws.Range("F9").Value = 18896
$ws.Range("F10").Value = 59
$ws.Range("F11").Value = 1992
$ws.Range("F12").Value = 687
$ws.Range("F14").Value = 350
$ws.Range("F15").Value = 619
$ws.Range("F16").Value = 202
$ws.Range("F17").Value = 213
$ws.Range("F19").Value = 329
$ws.Range("F20").Value = 48
$ws.Range("F21").Value = 214
$ws.Range("F23").Value = 122
$ws.Range("F24").Value = 5
$ws.Range("F25").Value = 15

# Sheet 2: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 19
$ws.Range("F5").Value = 176
$ws.Range("F9").Value = 112
$ws.Range("F10").Value = 237
$ws.Range("F18").Value = 23

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5913
$ws.Range("F3").Value = 589
$ws.Range("F4").Value = 564

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 589
$ws.Range("F4").Value = 564
$ws.Range("F6").Value = 19
$ws.Range("F7").Value = 169
$ws.Range("F9").Value = 1340
$ws.Range("F11").Value = 240
$ws.Range("F12").Value = 176
$ws.Range("F14").Value = 2555
$ws.Range("F15").Value = 942
$ws.Range("F16").Value = 18897
$ws.Range("F19").Value = 59
$ws.Range("F20").Value = 112
$ws.Range("F21").Value = 237
$ws.Range("F22").Value = 1992
$ws.Range("F23").Value = 687
$ws.Range("F25").Value = 350
$ws.Range("F26").Value = 619
$ws.Range("F27").Value = 202
$ws.Range("F28").Value = 213
$ws.Range("F32").Value = 329
$ws.Range("F33").Value = 48
$ws.Range("F35").Value = 214
$ws.Range("F38").Value = 122
$ws.Range("F39").Value = 23
$ws.Range("F41").Value = 5
$ws.Range("F44").Value = 15
